$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H97").Value = 200960
$ws.Range("J97").Value = 200960
$ws.Range("L97").Value = 602880
$ws.Range("N97").Value = -603872
$ws.Range("H101").Value = 393.33334
$ws.Range("I101").Value = 393.33334
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 1180.00002
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 441.9999800000001
$ws.Range("N101").ClearContents()
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H103").Value = 554.08826
$ws.Range("I103").Value = 411.8
$ws.Range("J103").Value = 578.62067
$ws.Range("K103").Value = 1235.4
$ws.Range("L103").Value = 1735.86201
$ws.Range("M103").Value = -649.4000000000001
$ws.Range("N103").Value = -2907.86201
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H106").Value = 4284.875
$ws.Range("I106").Value = 4284.875
$ws.Range("K106").Value = 4284.875
$ws.Range("M106").Value = -3653.875
$ws.Range("H121").Value = 1112.7142
$ws.Range("I121").Value = 2000
$ws.Range("J121").Value = 1044.4615
$ws.Range("K121").Value = 6000
$ws.Range("L121").Value = 3133.3845
$ws.Range("M121").Value = -4253
$ws.Range("N121").Value = -6627.3845
$ws.Range("H129").Value = 877.63336
$ws.Range("I129").Value = 348
$ws.Range("J129").Value = 959.11536
$ws.Range("K129").Value = 1044
$ws.Range("L129").Value = 2877.34608
$ws.Range("M129").Value = 3956
$ws.Range("N129").Value = -12877.34608

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5200.75
$ws.Range("I45").Value = 5023.25
$ws.Range("J45").Value = 5378.25
$ws.Range("K45").Value = 5023.25
$ws.Range("L45").Value = 5378.25
$ws.Range("M45").Value = -4646.25
$ws.Range("N45").Value = -6132.25
$ws.Range("H74").Value = 2913.739
$ws.Range("I74").Value = 2002.2858
$ws.Range("J74").Value = 4331.5557
$ws.Range("K74").Value = 2002.2858
$ws.Range("L74").Value = 4331.5557
$ws.Range("M74").Value = -1128.2858
$ws.Range("N74").Value = -6079.5557
$ws.Range("H77").Value = 2913.739
$ws.Range("I77").Value = 2002.2858
$ws.Range("J77").Value = 4331.5557
$ws.Range("K77").Value = 10011.429
$ws.Range("L77").Value = 21657.7785
$ws.Range("M77").Value = -5643.429
$ws.Range("N77").Value = -30393.7785

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 4583.4
$ws.Range("J93").Value = 4583.4
$ws.Range("L93").Value = 13750.2
$ws.Range("N93").Value = -17494.2
$ws.Range("H94").Value = 4333.3335
$ws.Range("J94").Value = 4000
$ws.Range("L94").Value = 12000
$ws.Range("N94").Value = -13352
$ws.Range("H95").Value = 4200
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 4200
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 12600
$ws.Range("M95").ClearContents()
$ws.Range("N95").Value = -16718
$ws.Range("H96").Value = 3800
$ws.Range("J96").Value = 3800
$ws.Range("L96").Value = 11400
$ws.Range("N96").Value = -15518
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H99").Value = 2425
$ws.Range("I99").Value = 2425
$ws.Range("K99").Value = 7275
$ws.Range("M99").Value = -5029
$ws.Range("H100").Value = 7311.8184
$ws.Range("I100").Value = 1180
$ws.Range("J100").Value = 7925
$ws.Range("K100").Value = 3540
$ws.Range("L100").Value = 23775
$ws.Range("M100").Value = -2729
$ws.Range("N100").Value = -25397
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H104").Value = 2749.25
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 2749.25
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 8247.75
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -13489.75
$ws.Range("H114").Value = 531.5
$ws.Range("I114").Value = 504.7
$ws.Range("J114").Value = 665.5
$ws.Range("K114").Value = 1514.1
$ws.Range("L114").Value = 1996.5
$ws.Range("M114").Value = 1739.9
$ws.Range("N114").Value = -8504.5
$ws.Range("H131").Value = 810.95
$ws.Range("J131").Value = 847.5806
$ws.Range("L131").Value = 2542.7418
$ws.Range("N131").Value = -12622.7418

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2578.7273
$ws.Range("I132").Value = 2062.2666
$ws.Range("J132").Value = 3685.4285
$ws.Range("K132").Value = 6186.7998
$ws.Range("L132").Value = 11056.2855
$ws.Range("M132").Value = -3656.7998
$ws.Range("N132").Value = -16116.2855

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1080
$ws.Range("I22").Value = 3433.3333
$ws.Range("J22").Value = 491.66666
$ws.Range("K22").Value = 3433.3333
$ws.Range("L22").Value = 491.66666
$ws.Range("M22").Value = -3138.3333
$ws.Range("N22").Value = -1081.66666
$ws.Range("H27").Value = 1080
$ws.Range("I27").Value = 3433.3333
$ws.Range("J27").Value = 491.66666
$ws.Range("K27").Value = 3433.3333
$ws.Range("L27").Value = 491.66666
$ws.Range("M27").Value = -3326.3333
$ws.Range("N27").Value = -705.66666
$ws.Range("H132").Value = 4581.591
$ws.Range("I132").Value = 5499.5
$ws.Range("J132").Value = 2133.8333
$ws.Range("K132").Value = 16498.5
$ws.Range("L132").Value = 6401.499899999999
$ws.Range("M132").Value = -13968.5
$ws.Range("N132").Value = -11461.4999
$ws.Range("H136").Value = 2087.8
$ws.Range("I136").Value = 1995.3182
$ws.Range("J136").Value = 2766
$ws.Range("K136").Value = 5985.9546
$ws.Range("L136").Value = 8298
$ws.Range("M136").Value = -3435.9546
$ws.Range("N136").Value = -13398
